$wb = $excel.ActiveWorkbook

# Sheet "展览" (1st tab) - F column ("想去人数") updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 22
$ws1.Range("F3").Value = 214
$ws1.Range("F7").Value = 106
$ws1.Range("F8").Value = 0
$ws1.Range("F10").Value = 724
$ws1.Range("F13").Value = 94
$ws1.Range("F14").Value = 257
$ws1.Range("F15").Value = 0
$ws1.Range("F17").Value = 135
$ws1.Range("F18").Value = 104
$ws1.Range("F20").Value = 0
$ws1.Range("F21").Value = 39
$ws1.Range("F22").Value = 0
$ws1.Range("F23").Value = 0
$ws1.Range("F26").Value = 0
$ws1.Range("F27").Value = 385
$ws1.Range("F28").Value = 0
$ws1.Range("F29").Value = 2520
$ws1.Range("F30").Value = 567
$ws1.Range("F31").Value = 0
$ws1.Range("F32").Value = 131
$ws1.Range("F33").Value = 252
$ws1.Range("F34").Value = 286
$ws1.Range("F36").Value = 150
$ws1.Range("F37").Value = 0
$ws1.Range("F38").Value = 928
$ws1.Range("F39").Value = 0
$ws1.Range("F40").Value = 44
$ws1.Range("F41").Value = 0
$ws1.Range("F42").Value = 0
$ws1.Range("F43").Value = 0
$ws1.Range("F44").Value = 72

# Sheet "全部类型" (4th tab) - F column ("想去人数") updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 22
$ws4.Range("F3").Value = 214
$ws4.Range("F4").Value = 0
$ws4.Range("F5").Value = 199
$ws4.Range("F7").Value = 106
$ws4.Range("F8").Value = 104
$ws4.Range("F9").Value = 98
$ws4.Range("F10").Value = 90
$ws4.Range("F15").Value = 257
$ws4.Range("F16").Value = 165
$ws4.Range("F18").Value = 135
$ws4.Range("F19").Value = 104
$ws4.Range("F21").Value = 0
$ws4.Range("F22").Value = 39
$ws4.Range("F23").Value = 36
$ws4.Range("F24").Value = 0
$ws4.Range("F25").Value = 531
$ws4.Range("F26").Value = 45
$ws4.Range("F27").Value = 0
$ws4.Range("F28").Value = 385
$ws4.Range("F29").Value = 0
$ws4.Range("F30").Value = 0
$ws4.Range("F34").Value = 252
$ws4.Range("F35").Value = 286
$ws4.Range("F37").Value = 0
$ws4.Range("F38").Value = 1545
$ws4.Range("F42").Value = 0
$ws4.Range("F43").Value = 473
$ws4.Range("F44").Value = 475
$ws4.Range("F45").Value = 0
